# Append 10 new device rows (new Mac-Addresses) to the
# master-reg_center_device_h table, continuing on from the last
# existing row (146 -> new rows 147-156).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$startRow = 147
$startMac = 3000166
$count = 10

for ($i = 0; $i -lt $count; $i++) {
    $r = $startRow + $i
    $ws.Cells.Item($r, 1).Value = 10001
    $ws.Cells.Item($r, 2).Value = $startMac + $i
    $ws.Cells.Item($r, 3).Value = "eng"
    $ws.Cells.Item($r, 4).Value = $true
    $ws.Cells.Item($r, 5).Value = "superadmin"
    $ws.Cells.Item($r, 6).Value = "now()"
    $ws.Cells.Item($r, 7).Value = "now()"
}

# Mirror the author's final view/selection state as closely as the
# runtime allows (scroll position isn't persisted by this headless
# engine outside of freeze-pane scenarios, but the active cell is).
$excel.ActiveWindow.ScrollRow = 140
$ws.Range("E155").Select()
